$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.21
$wsSummary.Range("B4").Value = 0.2
$wsSummary.Range("B5").Value = 0.06
$wsSummary.Range("B6").Value = 70
$wsSummary.Range("B7").Value = 32
$wsSummary.Range("B9").Value = 45.71

# --- Sheet: Strategy Status (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.21
$wsStatus.Range("D4").Value = 70
$wsStatus.Range("E4").Value = 0.2
$wsStatus.Range("F4").Value = 0.21
$wsStatus.Range("G4").Value = 45.71

# --- New trade row (#70) appended to both "All Trades" and "MarketMaking" sheets ---
$newRow = @{
    A = 70
    B = "2026-02-17"
    C = "12:54:16"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.92
    G = 0.95
    H = "CLOSED"
    I = 3.2609
    J = 0.03
    K = 100.21
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A71").Value = $newRow.A

    # B71/C71 hold date-like / time-like literal text ("2026-02-17",
    # "12:54:16"). A plain .Value assignment gets auto-coerced by the COM
    # layer into a date/time serial number, which is not what the source
    # workbook stores (plain text, t="inlineStr"/shared-string). Route the
    # literal text through a self-quoting formula, then collapse it back to
    # a literal value via copy/paste-values so no formula or number-format
    # style is left behind on the cell.
    $ws.Range("B71").Formula = '="' + $newRow.B + '"'
    $ws.Range("B71").Copy()
    $ws.Range("B71").PasteSpecial(-4163)
    $ws.Range("C71").Formula = '="' + $newRow.C + '"'
    $ws.Range("C71").Copy()
    $ws.Range("C71").PasteSpecial(-4163)

    $ws.Range("D71").Value = $newRow.D
    $ws.Range("E71").Value = $newRow.E
    $ws.Range("F71").Value = $newRow.F
    $ws.Range("G71").Value = $newRow.G
    $ws.Range("H71").Value = $newRow.H
    $ws.Range("I71").Value = $newRow.I
    $ws.Range("J71").Value = $newRow.J
    $ws.Range("K71").Value = $newRow.K
    $ws.Range("L71").Value = $newRow.L
    $ws.Range("M71").Value = $newRow.M
    $ws.Range("N71").Value = $newRow.N
    $ws.Range("O71").Value = $newRow.O
    $ws.Range("P71").Value = $newRow.P
    $ws.Range("Q71").Value = $newRow.Q
}
